# Update PLC data 2025-10-13 13:43:19
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 257
$ws.Range("C3").Value = 156049
$ws.Range("C4").Value = 147143
$ws.Range("C8").Value = 63.62
